$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Helper: write a numeric-looking value while keeping it stored as TEXT ---
# (Excel auto-converts a pure numeric string assigned through .Value into a
#  number; forcing the number format to Text first keeps it a string, then
#  the style is reset back to Normal so the cell ends up unstyled again.)
function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Row 2 — overall score (D2) and the new top "must have" / "good to have" matches
Set-TextValue $ws.Range("D2") "98.80"
$ws.Range("F2").Value = "docker : 3"
Set-TextValue $ws.Range("G2") "100.0"
$ws.Range("H2").Value = "go : 1"
Set-TextValue $ws.Range("I2") "50.0"

# Rows 3-7 — both CLUSTER MUST HAVE MATCH (F) and CLUSTER GOOD TO HAVE MATCH (H) populated
$ws.Range("F3").Value = "kafka : 3"
$ws.Range("H3").Value = "redshift : 2"

$ws.Range("F4").Value = "bi : 1"
$ws.Range("H4").Value = "dynamodb : 1"

$ws.Range("F5").Value = "aws : 4"
$ws.Range("H5").Value = "s3 : 1"

$ws.Range("F6").Value = "deployment : 1"
$ws.Range("H6").Value = "kinesis : 1"

$ws.Range("F7").Value = "kubernetes : 2"
$ws.Range("H7").Value = "spark : 1"

# Rows 8-11 — only CLUSTER MUST HAVE MATCH (F) populated
$ws.Range("F8").Value = "etl : 4"
$ws.Range("F9").Value = "mongodb : 1"
$ws.Range("F10").Value = "tools : 1"
$ws.Range("F11").Value = "amazon : 1"

# Rows 12-19 no longer have a CLUSTER MUST HAVE MATCH entry
$ws.Range("F12").ClearContents()
$ws.Range("F13").ClearContents()
$ws.Range("F14").ClearContents()
$ws.Range("F15").ClearContents()
$ws.Range("F16").ClearContents()
$ws.Range("F17").ClearContents()
$ws.Range("F18").ClearContents()
$ws.Range("F19").ClearContents()
